# Daily attendance processing - 2026-01-24 15:01:28
# Swap the order of names in the "Recorded By" column (column G) from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# wherever that exact value occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
